# Fruta / hortaliza, semanal
# Insert 3 new weekly observation rows (new rows 319-321) above the existing
# row that is currently 319 ("Feria Lagunitas de Puerto Montt" - Kiwi data).
# Excel shifts all the following rows (old 319-410) down to 322-413 and
# extends the used range to A1:T413 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("319:321").Insert()

# --- New row 319 ---
$ws.Range("A319").Value = 4
$ws.Range("B319").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C319").Value = "Los Lagos"
$ws.Range("D319").Value = 44876
$ws.Range("E319").Value = 10
$ws.Range("F319").Value = "Fruta"
$ws.Range("G319").Value = 100101
$ws.Range("H319").Value = "Berries"
$ws.Range("I319").Value = 100101007
$ws.Range("J319").Value = "Kiwi"
$ws.Range("K319").Value = "Hayward"
$ws.Range("L319").Value = "Especial"
$ws.Range("M319").Value = 200
$ws.Range("N319").Value = 20000
$ws.Range("O319").Value = 20000
$ws.Range("P319").Value = 20000
$ws.Range("Q319").Value = "`$/caja 15 kilos"
$ws.Range("R319").Value = "Región de O'Higgins"
$ws.Range("S319").Value = 1333
$ws.Range("T319").Value = 15

# --- New row 320 ---
$ws.Range("A320").Value = 4
$ws.Range("B320").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C320").Value = "Los Lagos"
$ws.Range("D320").Value = 44876
$ws.Range("E320").Value = 10
$ws.Range("F320").Value = "Fruta"
$ws.Range("G320").Value = 100101
$ws.Range("H320").Value = "Berries"
$ws.Range("I320").Value = 100101007
$ws.Range("J320").Value = "Kiwi"
$ws.Range("K320").Value = "Hayward"
$ws.Range("L320").Value = "Primera"
$ws.Range("M320").Value = 200
$ws.Range("N320").Value = 17000
$ws.Range("O320").Value = 17000
$ws.Range("P320").Value = 17000
$ws.Range("Q320").Value = "`$/caja 15 kilos"
$ws.Range("R320").Value = "Región de O'Higgins"
$ws.Range("S320").Value = 1133
$ws.Range("T320").Value = 15

# --- New row 321 ---
$ws.Range("A321").Value = 4
$ws.Range("B321").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C321").Value = "Los Lagos"
$ws.Range("D321").Value = 44876
$ws.Range("E321").Value = 10
$ws.Range("F321").Value = "Fruta"
$ws.Range("G321").Value = 100101
$ws.Range("H321").Value = "Berries"
$ws.Range("I321").Value = 100101007
$ws.Range("J321").Value = "Kiwi"
$ws.Range("K321").Value = "Hayward"
$ws.Range("L321").Value = "Segunda"
$ws.Range("M321").Value = 200
$ws.Range("N321").Value = 15000
$ws.Range("O321").Value = 15000
$ws.Range("P321").Value = 15000
$ws.Range("Q321").Value = "`$/caja 15 kilos"
$ws.Range("R321").Value = "Región de O'Higgins"
$ws.Range("S321").Value = 1000
$ws.Range("T321").Value = 15
